{"js": "// Office.js (Word JavaScript API) edit script.\n// Replaces the date line and each \"a\u00f7b=c, r\" division-table entry with its\n// updated value. Each (old, new) pair below is unique within the document,\n// so a plain text search + full-text replace is safe and unambiguous.\nconst pairs = [\n  [\"2024-09-02 Monday\", \"2024-09-03 Tuesday\"],\n  [\"29\u00f77=4, 1\", \"23\u00f72=11, 1\"],\n  [\"31\u00f77=4, 3\", \"24\u00f77=3, 3\"],\n  [\"52\u00f72=26, 0\", \"27\u00f73=9, 0\"],\n  [\"19\u00f77=2, 5\", \"67\u00f76=11, 1\"],\n  [\"13\u00f74=3, 1\", \"19\u00f72=9, 1\"],\n  [\"29\u00f73=9, 2\", \"75\u00f74=18, 3\"],\n  [\"63\u00f78=7, 7\", \"48\u00f75=9, 3\"],\n  [\"93\u00f78=11, 5\", \"10\u00f75=2, 0\"],\n  [\"60\u00f78=7, 4\", \"64\u00f77=9, 1\"],\n  [\"33\u00f75=6, 3\", \"61\u00f75=12, 1\"],\n  [\"62\u00f78=7, 6\", \"47\u00f79=5, 2\"],\n  [\"32\u00f76=5, 2\", \"94\u00f72=47, 0\"],\n  [\"23\u00f73=7, 2\", \"39\u00f77=5, 4\"],\n  [\"49\u00f73=16, 1\", \"51\u00f73=17, 0\"],\n  [\"57\u00f79=6, 3\", \"20\u00f72=10, 0\"],\n  [\"93\u00f79=10, 3\", \"99\u00f73=33, 0\"],\n  [\"33\u00f77=4, 5\", \"54\u00f75=10, 4\"],\n  [\"45\u00f74=11, 1\", \"68\u00f77=9, 5\"],\n  [\"18\u00f77=2, 4\", \"29\u00f78=3, 5\"],\n  [\"11\u00f79=1, 2\", \"78\u00f74=19, 2\"],\n  [\"50\u00f77=7, 1\", \"95\u00f72=47, 1\"],\n  [\"44\u00f73=14, 2\", \"93\u00f74=23, 1\"],\n  [\"23\u00f78=2, 7\", \"10\u00f76=1, 4\"],\n  [\"92\u00f76=15, 2\", \"77\u00f74=19, 1\"],\n  [\"77\u00f72=38, 1\", \"55\u00f77=7, 6\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, {matchCase: true, matchWholeWord: false});\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-09-02 Monday\", \"2024-09-03 Tuesday\"),\n    @(\"29\u00f77=4, 1\", \"23\u00f72=11, 1\"),\n    @(\"31\u00f77=4, 3\", \"24\u00f77=3, 3\"),\n    @(\"52\u00f72=26, 0\", \"27\u00f73=9, 0\"),\n    @(\"19\u00f77=2, 5\", \"67\u00f76=11, 1\"),\n    @(\"13\u00f74=3, 1\", \"19\u00f72=9, 1\"),\n    @(\"29\u00f73=9, 2\", \"75\u00f74=18, 3\"),\n    @(\"63\u00f78=7, 7\", \"48\u00f75=9, 3\"),\n    @(\"93\u00f78=11, 5\", \"10\u00f75=2, 0\"),\n    @(\"60\u00f78=7, 4\", \"64\u00f77=9, 1\"),\n    @(\"33\u00f75=6, 3\", \"61\u00f75=12, 1\"),\n    @(\"62\u00f78=7, 6\", \"47\u00f79=5, 2\"),\n    @(\"32\u00f76=5, 2\", \"94\u00f72=47, 0\"),\n    @(\"23\u00f73=7, 2\", \"39\u00f77=5, 4\"),\n    @(\"49\u00f73=16, 1\", \"51\u00f73=17, 0\"),\n    @(\"57\u00f79=6, 3\", \"20\u00f72=10, 0\"),\n    @(\"93\u00f79=10, 3\", \"99\u00f73=33, 0\"),\n    @(\"33\u00f77=4, 5\", \"54\u00f75=10, 4\"),\n    @(\"45\u00f74=11, 1\", \"68\u00f77=9, 5\"),\n    @(\"18\u00f77=2, 4\", \"29\u00f78=3, 5\"),\n    @(\"11\u00f79=1, 2\", \"78\u00f74=19, 2\"),\n    @(\"50\u00f77=7, 1\", \"95\u00f72=47, 1\"),\n    @(\"44\u00f73=14, 2\", \"93\u00f74=23, 1\"),\n    @(\"23\u00f78=2, 7\", \"10\u00f76=1, 4\"),\n    @(\"92\u00f76=15, 2\", \"77\u00f74=19, 1\"),\n    @(\"77\u00f72=38, 1\", \"55\u00f77=7, 6\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue=1, wdReplaceAll=2\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
